$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.023130727947523
$ws.Cells.Item(2, 4).Value = 1.035615442726687
$ws.Cells.Item(2, 5).Value = 1.023785036471742
$ws.Cells.Item(2, 6).Value = 1.021590676354961
$ws.Cells.Item(2, 9).Value = 1.034848642532224
$ws.Cells.Item(2, 10).Value = 1.02831267284185
$ws.Cells.Item(2, 11).Value = 1.038411539414883
$ws.Cells.Item(2, 12).Value = 1.026615455125279
$ws.Cells.Item(2, 13).Value = 1.024427556298778
$ws.Cells.Item(2, 14).Value = 1.029772994786401
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.024140357870903
$ws.Cells.Item(3, 4).Value = 1.036196161651538
$ws.Cells.Item(3, 5).Value = 1.024643008345294
$ws.Cells.Item(3, 6).Value = 1.02323862656823
$ws.Cells.Item(3, 9).Value = 1.035077221173152
$ws.Cells.Item(3, 10).Value = 1.028960337456806
$ws.Cells.Item(3, 11).Value = 1.038802297906844
$ws.Cells.Item(3, 12).Value = 1.027280240583787
$ws.Cells.Item(3, 13).Value = 1.025879690600248
$ws.Cells.Item(3, 14).Value = 1.030421579159396
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.024793234842583
$ws.Cells.Item(4, 4).Value = 1.036571360318531
$ws.Cells.Item(4, 5).Value = 1.0251981996443
$ws.Cells.Item(4, 6).Value = 1.024304414237874
$ws.Cells.Item(4, 9).Value = 1.035223376289552
$ws.Cells.Item(4, 10).Value = 1.029378440451319
$ws.Cells.Item(4, 11).Value = 1.039053880962959
$ws.Cells.Item(4, 12).Value = 1.027709791451731
$ws.Cells.Item(4, 13).Value = 1.026818320189108
$ws.Cells.Item(4, 14).Value = 1.030840275908117
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.025067604633605
$ws.Cells.Item(5, 4).Value = 1.036728957485635
$ws.Cells.Item(5, 5).Value = 1.025431608663841
$ws.Cells.Item(5, 6).Value = 1.024752347682966
$ws.Cells.Item(5, 9).Value = 1.035284400673741
$ws.Cells.Item(5, 10).Value = 1.029553977213777
$ws.Cells.Item(5, 11).Value = 1.03915934334674
$ws.Cells.Item(5, 12).Value = 1.027890229164885
$ws.Cells.Item(5, 13).Value = 1.027212687188839
$ws.Cells.Item(5, 14).Value = 1.031016061952913
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.025113666715553
$ws.Cells.Item(6, 4).Value = 1.036755410717863
$ws.Cells.Item(6, 5).Value = 1.025470799461902
$ws.Cells.Item(6, 6).Value = 1.024827550631849
$ws.Cells.Item(6, 9).Value = 1.035294622348442
$ws.Cells.Item(6, 10).Value = 1.029583436922744
$ws.Cells.Item(6, 11).Value = 1.039177033152315
$ws.Cells.Item(6, 12).Value = 1.027920516934989
$ws.Cells.Item(6, 13).Value = 1.027278889665241
$ws.Cells.Item(6, 14).Value = 1.031045563498045
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.024796901375482
$ws.Cells.Item(7, 4).Value = 1.036573466676672
$ws.Cells.Item(7, 5).Value = 1.025201318440719
$ws.Cells.Item(7, 6).Value = 1.02431040002372
$ws.Cells.Item(7, 9).Value = 1.035224193347831
$ws.Cells.Item(7, 10).Value = 1.02938078689969
$ws.Cells.Item(7, 11).Value = 1.039055291348343
$ws.Cells.Item(7, 12).Value = 1.02771220304028
$ws.Cells.Item(7, 13).Value = 1.026823590647719
$ws.Cells.Item(7, 14).Value = 1.030842625688714
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.023472024979165
$ws.Cells.Item(8, 4).Value = 1.035811815850306
$ws.Cells.Item(8, 5).Value = 1.024074986884039
$ws.Cells.Item(8, 6).Value = 1.022147726242332
$ws.Cells.Item(8, 9).Value = 1.034926254134058
$ws.Cells.Item(8, 10).Value = 1.028531757338986
$ws.Cells.Item(8, 11).Value = 1.038543859417651
$ws.Cells.Item(8, 12).Value = 1.026840249052092
$ws.Cells.Item(8, 13).Value = 1.02491852231425
$ws.Cells.Item(8, 14).Value = 1.029992390408649
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.021134147102561
$ws.Cells.Item(9, 4).Value = 1.034465404232668
$ws.Cells.Item(9, 5).Value = 1.022090430466267
$ws.Cells.Item(9, 6).Value = 1.018332316262657
$ws.Cells.Item(9, 9).Value = 1.034387845656116
$ws.Cells.Item(9, 10).Value = 1.027028125494133
$ws.Cells.Item(9, 11).Value = 1.037632996199009
$ws.Cells.Item(9, 12).Value = 1.02529906084883
$ws.Cells.Item(9, 13).Value = 1.021553626349326
$ws.Cells.Item(9, 14).Value = 1.028486623234107
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.019573284059292
$ws.Cells.Item(10, 4).Value = 1.033564972517911
$ws.Cells.Item(10, 5).Value = 1.020767490053333
$ws.Cells.Item(10, 6).Value = 1.015785235942943
$ws.Cells.Item(10, 9).Value = 1.034019901860458
$ws.Cells.Item(10, 10).Value = 1.02602059064839
$ws.Cells.Item(10, 11).Value = 1.037019293778018
$ws.Cells.Item(10, 12).Value = 1.024268410245704
$ws.Cells.Item(10, 13).Value = 1.01930464912694
$ws.Cells.Item(10, 14).Value = 1.027477657573318
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.018896855186722
$ws.Cells.Item(11, 4).Value = 1.033174417317162
$ws.Cells.Item(11, 5).Value = 1.020194658413008
$ws.Cells.Item(11, 6).Value = 1.014681398391499
$ws.Cells.Item(11, 9).Value = 1.033858442865716
$ws.Cells.Item(11, 10).Value = 1.025583092386188
$ws.Cells.Item(11, 11).Value = 1.036752027482242
$ws.Cells.Item(11, 12).Value = 1.023821361260666
$ws.Cells.Item(11, 13).Value = 1.018329375512949
$ws.Cells.Item(11, 14).Value = 1.0270395380134
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.018645512881558
$ws.Cells.Item(12, 4).Value = 1.033029249077126
$ws.Cells.Item(12, 5).Value = 1.019981884081216
$ws.Cells.Item(12, 6).Value = 1.014271234852776
$ws.Cells.Item(12, 9).Value = 1.033798148719386
$ws.Cells.Item(12, 10).Value = 1.025420400300822
$ws.Cells.Item(12, 11).Value = 1.036652523549112
$ws.Cells.Item(12, 12).Value = 1.023655190750662
$ws.Cells.Item(12, 13).Value = 1.017966889848684
$ws.Cells.Item(12, 14).Value = 1.026876614886614
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.018699430600168
$ws.Cells.Item(13, 4).Value = 1.03306039260786
$ws.Cells.Item(13, 5).Value = 1.020027524870858
$ws.Cells.Item(13, 6).Value = 1.014359223211858
$ws.Cells.Item(13, 9).Value = 1.03381109654905
$ws.Cells.Item(13, 10).Value = 1.025455306728918
$ws.Cells.Item(13, 11).Value = 1.036673877851972
$ws.Cells.Item(13, 12).Value = 1.023690840179992
$ws.Cells.Item(13, 13).Value = 1.018044654585288
$ws.Cells.Item(13, 14).Value = 1.026911570885841
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.018876080927553
$ws.Cells.Item(14, 4).Value = 1.03316241967266
$ws.Cells.Item(14, 5).Value = 1.020177070407245
$ws.Cells.Item(14, 6).Value = 1.014647497237771
$ws.Cells.Item(14, 9).Value = 1.033853465483113
$ws.Cells.Item(14, 10).Value = 1.025569648003366
$ws.Cells.Item(14, 11).Value = 1.036743807132107
$ws.Cells.Item(14, 12).Value = 1.023807627938333
$ws.Cells.Item(14, 13).Value = 1.018299416981152
$ws.Cells.Item(14, 14).Value = 1.027026074538013
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.018984909512584
$ws.Cells.Item(15, 4).Value = 1.033225268868046
$ws.Cells.Item(15, 5).Value = 1.020269210438541
$ws.Cells.Item(15, 6).Value = 1.014825092359209
$ws.Cells.Item(15, 9).Value = 1.033879527833253
$ws.Cells.Item(15, 10).Value = 1.02564007278918
$ws.Cells.Item(15, 11).Value = 1.036786862492709
$ws.Cells.Item(15, 12).Value = 1.02387956925261
$ws.Cells.Item(15, 13).Value = 1.018456354378644
$ws.Cells.Item(15, 14).Value = 1.027096599335101
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.019618164375098
$ws.Cells.Item(16, 4).Value = 1.033590878499464
$ws.Cells.Item(16, 5).Value = 1.020805507160016
$ws.Cells.Item(16, 6).Value = 1.015858473569552
$ws.Cells.Item(16, 9).Value = 1.034030572336776
$ws.Cells.Item(16, 10).Value = 1.026049599979138
$ws.Cells.Item(16, 11).Value = 1.037036999179838
$ws.Cells.Item(16, 12).Value = 1.024298063108689
$ws.Cells.Item(16, 13).Value = 1.019369343598888
$ws.Cells.Item(16, 14).Value = 1.027506708100643
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.0200152360363
$ws.Cells.Item(17, 4).Value = 1.033820039049194
$ws.Cells.Item(17, 5).Value = 1.021141914189604
$ws.Cells.Item(17, 6).Value = 1.016506429892374
$ws.Cells.Item(17, 9).Value = 1.03412474642748
$ws.Cells.Item(17, 10).Value = 1.026306155746577
$ws.Cells.Item(17, 11).Value = 1.03719349421784
$ws.Cells.Item(17, 12).Value = 1.024560366370299
$ws.Cells.Item(17, 13).Value = 1.019941643673488
$ws.Cells.Item(17, 14).Value = 1.027763628206698
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.020246786826893
$ws.Cells.Item(18, 4).Value = 1.033953640559201
$ws.Cells.Item(18, 5).Value = 1.021338135758995
$ws.Cells.Item(18, 6).Value = 1.016884282225702
$ws.Cells.Item(18, 9).Value = 1.034179470404364
$ws.Cells.Item(18, 10).Value = 1.026455681876889
$ws.Cells.Item(18, 11).Value = 1.037284627564007
$ws.Cells.Item(18, 12).Value = 1.024713289143813
$ws.Cells.Item(18, 13).Value = 1.020275317069991
$ws.Cells.Item(18, 14).Value = 1.027913366681267
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.02032573048652
$ws.Cells.Item(19, 4).Value = 1.033999184346213
$ws.Cells.Item(19, 5).Value = 1.021405042440925
$ws.Cells.Item(19, 6).Value = 1.017013105167007
$ws.Cells.Item(19, 9).Value = 1.034198094896226
$ws.Cells.Item(19, 10).Value = 1.026506646355532
$ws.Cells.Item(19, 11).Value = 1.037315676638593
$ws.Cells.Item(19, 12).Value = 1.024765419284198
$ws.Cells.Item(19, 13).Value = 1.020389067632792
$ws.Cells.Item(19, 14).Value = 1.027964403535315
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.019972639642402
$ws.Cells.Item(20, 4).Value = 1.033795458900501
$ws.Cells.Item(20, 5).Value = 1.021105820796267
$ws.Cells.Item(20, 6).Value = 1.016436919643255
$ws.Cells.Item(20, 9).Value = 1.034114663764845
$ws.Cells.Item(20, 10).Value = 1.026278642006359
$ws.Cells.Item(20, 11).Value = 1.037176719041342
$ws.Cells.Item(20, 12).Value = 1.024532231401976
$ws.Cells.Item(20, 13).Value = 1.019880255779969
$ws.Cells.Item(20, 14).Value = 1.027736075393814
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.018824064210195
$ws.Cells.Item(21, 4).Value = 1.033132377961415
$ws.Cells.Item(21, 5).Value = 1.020133032956491
$ws.Cells.Item(21, 6).Value = 1.014562611908499
$ws.Cells.Item(21, 9).Value = 1.033840997751988
$ws.Cells.Item(21, 10).Value = 1.025535982484503
$ws.Cells.Item(21, 11).Value = 1.03672322102956
$ws.Cells.Item(21, 12).Value = 1.023773240080379
$ws.Cells.Item(21, 13).Value = 1.018224402094255
$ws.Cells.Item(21, 14).Value = 1.026992361210253
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.018101407941874
$ws.Cells.Item(22, 4).Value = 1.032714902318115
$ws.Cells.Item(22, 5).Value = 1.019521407160858
$ws.Cells.Item(22, 6).Value = 1.013383293136726
$ws.Cells.Item(22, 9).Value = 1.033667075332499
$ws.Cells.Item(22, 10).Value = 1.025067967710695
$ws.Cells.Item(22, 11).Value = 1.036436762181232
$ws.Cells.Item(22, 12).Value = 1.023295356995964
$ws.Cells.Item(22, 13).Value = 1.017181992033012
$ws.Cells.Item(22, 14).Value = 1.026523681801784
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.018484549772249
$ws.Cells.Item(23, 4).Value = 1.032936267834741
$ws.Cells.Item(23, 5).Value = 1.019845641376487
$ws.Cells.Item(23, 6).Value = 1.014008557427935
$ws.Cells.Item(23, 9).Value = 1.033759451005741
$ws.Cells.Item(23, 10).Value = 1.025316173492014
$ws.Cells.Item(23, 11).Value = 1.036588745057832
$ws.Cells.Item(23, 12).Value = 1.023548756144159
$ws.Cells.Item(23, 13).Value = 1.017734719946457
$ws.Cells.Item(23, 14).Value = 1.026772240063781
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.019991887286834
$ws.Cells.Item(24, 4).Value = 1.033806565809301
$ws.Cells.Item(24, 5).Value = 1.021122129843588
$ws.Cells.Item(24, 6).Value = 1.016468328607724
$ws.Cells.Item(24, 9).Value = 1.034119220323033
$ws.Cells.Item(24, 10).Value = 1.026291074647058
$ws.Cells.Item(24, 11).Value = 1.037184299477417
$ws.Cells.Item(24, 12).Value = 1.024544944612316
$ws.Cells.Item(24, 13).Value = 1.019907994756065
$ws.Cells.Item(24, 14).Value = 1.027748525690288
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.021738940503268
$ws.Cells.Item(25, 4).Value = 1.034813985797939
$ws.Cells.Item(25, 5).Value = 1.022603466754936
$ws.Cells.Item(25, 6).Value = 1.019319270583836
$ws.Cells.Item(25, 9).Value = 1.034528624690542
$ws.Cells.Item(25, 10).Value = 1.027417748059823
$ws.Cells.Item(25, 11).Value = 1.037869617071563
$ws.Cells.Item(25, 12).Value = 1.025698054837784
$ws.Cells.Item(25, 13).Value = 1.022424511126915
$ws.Cells.Item(25, 14).Value = 1.028876799108531
